# Update "想去人数" (want-to-go count) values in column F across the
# workbook's sheets, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 1747
$ws.Range("F5").Value  = 438
$ws.Range("F7").Value  = 61
$ws.Range("F9").Value  = 292
$ws.Range("F10").Value = 1675
$ws.Range("F11").Value = 332
$ws.Range("F12").Value = 1384
$ws.Range("F13").Value = 784
$ws.Range("F16").Value = 12604
$ws.Range("F17").Value = 12636
$ws.Range("F19").Value = 731
$ws.Range("F21").Value = 296
$ws.Range("F23").Value = 488
$ws.Range("F24").Value = 1975
$ws.Range("F27").Value = 227
$ws.Range("F28").Value = 660

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 48
$ws.Range("F6").Value = 10
$ws.Range("F7").Value = 5

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 149

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 149
$ws.Range("F6").Value  = 1747
$ws.Range("F7").Value  = 438
$ws.Range("F10").Value = 61
$ws.Range("F13").Value = 48
$ws.Range("F14").Value = 292
$ws.Range("F15").Value = 1675
$ws.Range("F16").Value = 332
$ws.Range("F17").Value = 1384
$ws.Range("F18").Value = 784
$ws.Range("F22").Value = 12604
$ws.Range("F23").Value = 12636
$ws.Range("F25").Value = 731
$ws.Range("F27").Value = 296
$ws.Range("F29").Value = 488
$ws.Range("F30").Value = 10
$ws.Range("F31").Value = 5
$ws.Range("F32").Value = 1975
$ws.Range("F37").Value = 227
$ws.Range("F38").Value = 660
